$d = $word.ActiveDocument
$d.Content.Find.Execute("authors who create or update conroed Word documents in the", $true, $false, $false, $false, $false,
                         $true, 1, $false, "authors who create or update controlled Word documents in the", 2)
$d.Content.Find.Execute(" folder of the GitHub repository;", $true, $false, $false, $false, $false,
                         $true, 1, $false, " folder of the GitHub repository;", 2)
